# Apply the "update new orleans xlsx files" edit:
#   1. hotel_info gains a new "State" column (value "Louisiana") inserted
#      between "Hotel_Name" and "City".
#   2. The sheet order is swapped so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Insert a new column C ("State") in hotel_info, shifting City (and the
# rest of the header row / data row) one column to the right.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# Move review_info to be the first sheet (hotel_info becomes second).
$wsReview.Move($wb.Worksheets.Item(1))
